# Insert two new data rows (710 and 711) into the "Ajo" (garlic) price sheet,
# pushing the former rows 710..765 down to become rows 712..767.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two whole rows at position 710 (shifts everything below down by 2).
$ws.Range("A710:A711").EntireRow.Insert()

# --- New row 710 ---
$ws.Cells.Item(710, 1).Value = 6
$ws.Cells.Item(710, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(710, 3).Value = "Metropolitana"
$ws.Cells.Item(710, 4).Value = 44769
$ws.Cells.Item(710, 5).Value = 13
$ws.Cells.Item(710, 6).Value = 100112003
$ws.Cells.Item(710, 7).Value = "Ajo"
$ws.Cells.Item(710, 8).Value = "Rosado"
$ws.Cells.Item(710, 9).Value = "1a (guarda)"
$ws.Cells.Item(710, 10).Value = 300
$ws.Cells.Item(710, 11).Value = 24000
$ws.Cells.Item(710, 12).Value = 24000
$ws.Cells.Item(710, 13).Value = 24000
$ws.Cells.Item(710, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(710, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(710, 16).Value = 2400
$ws.Cells.Item(710, 17).Value = 10
$ws.Cells.Item(710, 18).Value = "Hortaliza"

# --- New row 711 ---
$ws.Cells.Item(711, 1).Value = 6
$ws.Cells.Item(711, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(711, 3).Value = "Metropolitana"
$ws.Cells.Item(711, 4).Value = 44769
$ws.Cells.Item(711, 5).Value = 13
$ws.Cells.Item(711, 6).Value = 100112003
$ws.Cells.Item(711, 7).Value = "Ajo"
$ws.Cells.Item(711, 8).Value = "Rosado"
$ws.Cells.Item(711, 9).Value = "2a (guarda)"
$ws.Cells.Item(711, 10).Value = 200
$ws.Cells.Item(711, 11).Value = 18000
$ws.Cells.Item(711, 12).Value = 18000
$ws.Cells.Item(711, 13).Value = 18000
$ws.Cells.Item(711, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(711, 15).Value = "Provincia de Talagante"
$ws.Cells.Item(711, 16).Value = 1800
$ws.Cells.Item(711, 17).Value = 10
$ws.Cells.Item(711, 18).Value = "Hortaliza"
